$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("216:219").Insert()
Write-Host $ws.Range("D216").Value2
Write-Host $ws.Range("D220").Value2
Write-Host $ws.Range("K220").Text
